$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 2: merge the two existing runs into a single run (text itself
#    is unchanged once concatenated).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Con riferimento al database IFTS scrivi la query che trova tutte le regioni che appartengono all’area geografica Sud",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Con riferimento al database IFTS scrivi la query che trova tutte le regioni che appartengono all’area geografica Sud",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Paragraph 3: merge the two existing runs into a single run (text itself
#    is unchanged once concatenated). The trailing bookmark currently
#    attached to this paragraph will be relocated later.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Con riferimento al database IFTS scrivi la query che trova tutte le prenotazioni che hanno importo superiore a 200",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Con riferimento al database IFTS scrivi la query che trova tutte le prenotazioni che hanno importo superiore a 200",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Detach the "_GoBack" bookmark from paragraph 3; it will be re-created at
#    the end of the new, final paragraph once that paragraph exists.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 4) Insert three new list paragraphs after paragraph 3, inheriting the
#    "Paragrafoelenco" / numId 8 list formatting already used by the
#    surrounding items.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$afterP3 = $p3.Range.InsertParagraphAfter()

$newPara1 = $d.Paragraphs(4)
$newPara1.Range.InsertParagraphAfter() | Out-Null

$newPara2 = $d.Paragraphs(5)
$newPara2.Range.InsertParagraphAfter() | Out-Null

# Paragraph 4 (single run).
$xmlPara4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Con riferimento al database IFTS scrivi la query che trova il saldo (importo – caparra) di tutte le prenotazioni</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(4).Range.InsertXML($xmlPara4)

# Paragraph 5 (four runs).
$xmlPara5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Con riferimento al database IFTS scrivi la query che trova il saldo (importo – caparra) d</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">elle </w:t></w:r>' +
    '<w:r><w:t>prenotazioni</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> per le strutture 2 stelle e 3 stelle</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(5).Range.InsertXML($xmlPara5)

# Paragraph 6 (five runs); bookmark is re-added after this text is in place.
$xmlPara6 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Con riferimento al database IFTS scrivi la query che trova tutte le prenotazioni che hanno </w:t></w:r>' +
    '<w:r><w:t>saldo (importo – caparra)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> superiore a </w:t></w:r>' +
    '<w:r><w:t>1</w:t></w:r>' +
    '<w:r><w:t>00</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(6).Range.InsertXML($xmlPara6)

# ---------------------------------------------------------------------------
# 5) Re-attach the "_GoBack" bookmark at the very end of the new paragraph 6
#    (right before its paragraph mark), matching the original placement.
# ---------------------------------------------------------------------------
$finalPara = $d.Paragraphs(6)
$bmPoint = $d.Range($finalPara.Range.End, $finalPara.Range.End)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# ---------------------------------------------------------------------------
# 6) Remove the trailing empty "Paragrafoelenco" paragraph left over at the
#    end of the document body.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Delete() | Out-Null
